$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values look numeric,
# so Excel keeps them as text (matching the original inline-string cell type).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates row by row
# Row 2
$ws.Range("D2").Value = "69.433.06"
$ws.Range("E2").Value = "  -1.58%  "
# Row 3
$ws.Range("D3").Value = "3.533.91"
$ws.Range("E3").Value = "  -2.68%  "
# Row 4
$ws.Range("E4").Value = "  +0.30%  "
# Row 5
$ws.Range("D5").Value = "585.80"
$ws.Range("E5").Value = "  -3.02%  "
# Row 6
$ws.Range("D6").Value = "193.93"
$ws.Range("E6").Value = "  -1.41%  "
# Row 7
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  -3.19%  "
# Row 8
$ws.Range("E8").Value = "  +0.09%  "
# Row 9
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").Value = "  -1.26%  "
# Row 10
$ws.Range("D10").Value = "0.621"
$ws.Range("E10").Value = "  -4.15%  "
# Row 11
$ws.Range("D11").Value = "52.36"
$ws.Range("E11").Value = "  -1.95%  "
# Row 12
$ws.Range("D12").Value = "0.0000287"
$ws.Range("E12").Value = "  -5.35%  "
# Row 13
$ws.Range("D13").Value = "9.21"
$ws.Range("E13").Value = "  -3.80%  "
# Row 14
$ws.Range("D14").Value = "4.106.84"
$ws.Range("E14").Value = "  -2.19%  "
# Row 15
$ws.Range("D15").Value = "648.07"
$ws.Range("E15").Value = "  +7.90%  "
# Row 16
$ws.Range("D16").Value = "69.573.26"
$ws.Range("E16").Value = "  -1.40%  "
# Row 17
$ws.Range("D17").Value = "3.539.02"
$ws.Range("E17").Value = "  -1.83%  "
# Row 18
$ws.Range("D18").Value = "12.53"
$ws.Range("E18").Value = "  -3.82%  "
# Row 19
$ws.Range("E19").Value = "  -1.50%  "
# Row 20
$ws.Range("D20").Value = "18.34"
$ws.Range("E20").Value = "  -3.73%  "
# Row 21
$ws.Range("D21").Value = "0.957"
$ws.Range("E21").Value = "  -4.00%  "
# Row 22
$ws.Range("D22").Value = "18.09"
$ws.Range("E22").Value = "  -0.67%  "
# Row 23
$ws.Range("D23").Value = "5.49"
$ws.Range("E23").Value = "  +5.70%  "
# Row 24
$ws.Range("D24").Value = "101.94"
$ws.Range("E24").Value = "  -0.05%  "
# Row 25
$ws.Range("D25").Value = "4.37"
$ws.Range("E25").Value = "  -5.59%  "
# Row 26
$ws.Range("D26").Value = "2.91"
$ws.Range("E26").Value = "  -3.29%  "
# Row 27
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -5.42%  "
# Row 28
$ws.Range("D28").Value = "9.51"
$ws.Range("E28").Value = "  -1.87%  "
# Row 29
$ws.Range("D29").Value = "33.01"
$ws.Range("E29").Value = "  -2.55%  "
# Row 30
$ws.Range("D30").Value = "4.10"
$ws.Range("E30").Value = "  -11.58%  "
# Row 31
$ws.Range("D31").Value = "6.75"
$ws.Range("E31").Value = "  -7.24%  "
# Row 32
$ws.Range("D32").Value = "11.69"
$ws.Range("E32").Value = "  -4.95%  "
# Row 33
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  -7.56%  "
# Row 34
$ws.Range("D34").Value = "61.59"
$ws.Range("E34").Value = "  -2.98%  "
# Row 35
$ws.Range("D35").Value = "3.725.46"
$ws.Range("E35").Value = "  -4.95%  "
# Row 36
$ws.Range("E36").Value = "  -0.13%  "
# Row 37
$ws.Range("D37").Value = "0.0₃0803"
$ws.Range("E37").Value = "  -9.97%  "
# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "511.08"
$ws.Range("E38").Value = "  -5.17%  "
# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.59"
$ws.Range("E39").Value = "  +1.39%  "
# Row 40
$ws.Range("D40").Value = "2.94"
$ws.Range("E40").Value = "  -4.45%  "
# Row 41
$ws.Range("D41").Value = "0.369"
$ws.Range("E41").Value = "  -5.42%  "
# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").Value = "  +0.18%  "
# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "34.62"
$ws.Range("E43").Value = "  -6.51%  "
# Row 44
$ws.Range("D44").Value = "0.0446"
$ws.Range("E44").Value = "  -3.66%  "
# Row 45
$ws.Range("D45").Value = "3.42"
$ws.Range("E45").Value = "  +1.20%  "
# Row 46
$ws.Range("D46").Value = "2.85"
$ws.Range("E46").Value = "  -1.16%  "
# Row 47
$ws.Range("E47").Value = "  -3.39%  "
# Row 48
$ws.Range("E48").Value = "  -0.06%  "
# Row 49
$ws.Range("D49").Value = "8.23"
$ws.Range("E49").Value = "  -4.16%  "
# Row 50
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "2.71"
$ws.Range("E50").Value = "  +57.67%  "
# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.154"
$ws.Range("E51").Value = "  +1.76%  "

Write-Output "Applied cryptos update"
